$d = $word.ActiveDocument

# 1. Merge the split "Mon Sep 17" / " 10:57:12 PDT 2017" runs into a single run.
$d.Content.Find.Execute("Mon Sep 17 10:57:12 PDT 2017", $true, $false, $false, $false, $false, $true, 1, $false, "Mon Sep 17 10:57:12 PDT 2017", 2) | Out-Null

# 2. Remove the two blank "PlainText" paragraphs that follow the
#    "Amount balance ... - 11244.0" line (entry for 17/09/2017).
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("- 11244.0")) {
        $d.Paragraphs.Item($i + 1).Range.Delete()
        $d.Paragraphs.Item($i + 1).Range.Delete()
        break
    }
}
